# "New Environment Documentation Update"
# Insert a new "Building" entry into the Static/Environment section of the
# Asset List, and bump the Player (gun) model's dimensions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room under "Ramp Platform" (row 7) for the new Building row, plus
# five blank spacer rows, so the "Characters" section header lands back on
# row 14 (it was previously on row 8).
$ws.Range("B8:B13").EntireRow.Insert()

# The inserted rows copy row 7's fill/format; strip that back out so rows
# 9-13 are genuinely blank and row 8 only carries the two new cells.
$ws.Range("B9:G13").Clear()
$ws.Range("D8:G8").Clear()

# New Static/Environment row: a residential building asset.
$ws.Range("B8").Value = "Building"
$ws.Range("C8").Value = "City residential building, slightly taller than the platform"

# Player (gun) dimensions changed (now at row 18 after the insert above).
$ws.Range("D18").Value = "0.5 x 0.5 x 0.5 Cube"

# Match the author's final selection.
$ws.Range("D8").Select() | Out-Null
